# ---------------------------------------------------------------------------
# Refresh the "Price" / "Volume(1h)" columns of the cryptos sheet to match the
# newer coinranking.com snapshot pulled in by the scheduled GitHub Actions run
# (also re-sorts the Toncoin / InjectiveProtocol pair - rows 35 & 36 swap).
#
# All of these source cells are plain text (t="inlineStr"/shared-string "Text"
# cells holding things like "51.007.38" or "  -1.12%  "), never real numbers -
# the dotted price strings are thousand-separated, not decimals. A handful of
# the new values (e.g. "0.590", "16.50") are themselves valid decimal literals,
# so a naive Range.Value assignment would let Excel "helpfully" reinterpret them
# as numbers (silently dropping the trailing zero / introducing float noise).
# For those cells we briefly force text storage via a leading apostrophe and
# then restore the cell to the default "Normal" style so no formatting is left
# behind.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.007.38"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.945.72"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'375.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "'101.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "'36.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'0.0850"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "3.402.87"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "'7.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'11.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +52.55%  "
$ws.Range("D17").Value = "2.943.34"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "'0.998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "50.967.43"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -5.94%  "
$ws.Range("D21").Value = "'12.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'266.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "'68.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.42%  "
$ws.Range("D26").Value = "'8.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").Value = "'7.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'25.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("D33").Value = "'9.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'33.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("D37").Value = "'0.0442"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'3.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Value = "'16.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "'2.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("D44").Value = "'120.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "'21.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").Value = "'3.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").Value = "1.992.33"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("E51").Value = "  -1.36%  "
